$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.23
$ws.Range("P2").Value = 2.5
$ws.Range("R2").Value = 1.6
$ws.Range("S2").Value = 2.34
$ws.Range("T2").Value = 1.58
$ws.Range("U2").Value = 2.4
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 8
$ws.Range("G3").Value = 1.46
$ws.Range("H3").Value = 8.6
$ws.Range("I3").Value = 8.800000000000001
$ws.Range("J3").Value = 5.1
$ws.Range("K3").Value = 5.2
$ws.Range("N3").Value = 5.5
$ws.Range("T3").Value = 1.87
$ws.Range("U3").Value = 2.08
$ws.Range("W3").Value = 3.2
$ws.Range("G4").Value = 2.06
$ws.Range("I4").Value = 3.7
$ws.Range("N4").Value = 5.2
$ws.Range("P4").Value = 2.44
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.56
$ws.Range("V4").Value = 1.37
$ws.Range("G5").Value = 2.1
$ws.Range("J5").Value = 3.7
$ws.Range("L5").Value = 1.2
$ws.Range("P5").Value = 2.72
$ws.Range("Q5").Value = 1.39
$ws.Range("R5").Value = 1.81
$ws.Range("S5").Value = 1.96
$ws.Range("T5").Value = 1.45
$ws.Range("U5").Value = 2.74
$ws.Range("W5").Value = 1.91
$ws.Range("AC5").Value = 13
$ws.Range("AH5").Value = 18.5
$ws.Range("F6").Value = 2.7
$ws.Range("G6").Value = 3.05
$ws.Range("H6").Value = 2.38
$ws.Range("I6").Value = 2.66
$ws.Range("J6").Value = 3.75
$ws.Range("K6").Value = 4.4
$ws.Range("L6").Value = 1.25
$ws.Range("N6").Value = 5.1
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 2.4
$ws.Range("Q6").Value = 1.59
$ws.Range("R6").Value = 1.57
$ws.Range("S6").Value = 2.46
$ws.Range("T6").Value = 1.54
$ws.Range("U6").Value = 2.5
$ws.Range("V6").Value = 1.6
$ws.Range("W6").Value = 1.5
$ws.Range("Y6").Value = 15.5
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 25
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 15.5
$ws.Range("AN6").Value = 18
$ws.Range("AO6").Value = 14.5
$ws.Range("F7").Value = 1.49
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 6.6
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 4.3
$ws.Range("L7").Value = 1.28
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 4.2
$ws.Range("R7").Value = 1.44
$ws.Range("S7").Value = 2.54
$ws.Range("T7").Value = 1.84
$ws.Range("U7").Value = 1.94
$ws.Range("V7").Value = 1.12
$ws.Range("W7").Value = 2.62
$ws.Range("X7").Value = 24
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 75
$ws.Range("AB7").Value = 11
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 34
$ws.Range("AF7").Value = 12
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 24
$ws.Range("AJ7").Value = 17
$ws.Range("AK7").Value = 19.5
$ws.Range("AL7").Value = 42
$ws.Range("AN7").Value = 8.800000000000001
$ws.Range("L8").Value = 1.45
$ws.Range("O8").Value = 1.38
$ws.Range("Q8").Value = 2.16
$ws.Range("F9").Value = 2.62
$ws.Range("I9").Value = 2.66
$ws.Range("L9").Value = 1.26
$ws.Range("O9").Value = 1.17
$ws.Range("V9").Value = 1.6
$ws.Range("F10").Value = 9.6
$ws.Range("K10").Value = 5.7
$ws.Range("R10").Value = 1.6
$ws.Range("S10").Value = 2.58
$ws.Range("T10").Value = 1.95
$ws.Range("U10").Value = 1.98
$ws.Range("AA10").Value = 11.5
$ws.Range("P11").Value = 2.96
$ws.Range("U11").Value = 1.69
$ws.Range("AE11").Value = 450
$ws.Range("AJ11").Value = 8.199999999999999
$ws.Range("F12").Value = 1.3
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 1.2
$ws.Range("T12").Value = 1.72
$ws.Range("Z12").Value = 120
$ws.Range("AE12").Value = 130
$ws.Range("AI12").Value = 95
$ws.Range("AL12").Value = 25
$ws.Range("AO12").Value = 95
$ws.Range("G13").Value = 6.4
$ws.Range("H13").Value = 1.64
$ws.Range("I13").Value = 1.65
$ws.Range("J13").Value = 4.3
$ws.Range("K13").Value = 4.4
$ws.Range("Q13").Value = 1.83
$ws.Range("V13").Value = 2.54
$ws.Range("AM13").Value = 110
$ws.Range("L14").Value = 1.34
$ws.Range("T14").Value = 1.63
$ws.Range("F15").Value = 2.76
$ws.Range("G15").Value = 3.4
$ws.Range("H15").Value = 2.18
$ws.Range("I15").Value = 2.58
$ws.Range("J15").Value = 3.1
$ws.Range("K15").Value = 4.5
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 4
$ws.Range("P15").Value = 2.22
$ws.Range("Q15").Value = 1.52
$ws.Range("R15").Value = 1.53
$ws.Range("S15").Value = 2.32
$ws.Range("T15").Value = 1.54
$ws.Range("U15").Value = 2.36
$ws.Range("V15").Value = 1.63
$ws.Range("W15").Value = 1.41
$ws.Range("X15").Value = 27
$ws.Range("Y15").Value = 16
$ws.Range("Z15").Value = 22
$ws.Range("AA15").Value = 38
$ws.Range("AB15").Value = 20
$ws.Range("AC15").Value = 11.5
$ws.Range("AD15").Value = 13.5
$ws.Range("AE15").Value = 28
$ws.Range("AF15").Value = 29
$ws.Range("AG15").Value = 16.5
$ws.Range("AH15").Value = 18.5
$ws.Range("AI15").Value = 34
$ws.Range("AJ15").Value = 60
$ws.Range("AK15").Value = 38
$ws.Range("AL15").Value = 44
$ws.Range("AN15").Value = 25
$ws.Range("AO15").Value = 16
$ws.Range("F16").Value = 2.5
$ws.Range("L16").Value = 1.36
$ws.Range("M16").Value = 1.08
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 1.31
$ws.Range("U16").Value = 2.06
$ws.Range("V16").Value = 1.45
$ws.Range("AB16").Value = 12.5
$ws.Range("AD16").Value = 16
$ws.Range("AN16").Value = 32
